# Apply the edit described by the diff:
# 1. Rename the worksheet from "Sheet1" to "books"
# 2. Append 20 new book rows (rows 7-26) with Title, Authors, Genres, Price
# 3. Update selection to E34 (next empty row after the data)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "books"

# New book rows to append starting at row 7
$rows = @(
    @("The Sisters Brothers", "Patrick DeWitt", "Western;Literature;Historical", 8.99),
    @("The Essex Serpent", "Sarah Perry", "Literature;Historical", 8.99),
    @("Why I’m No Longer Talking to White People About Race", "Reni Eddo-Lodge", "Non-Fiction;Social Sciences", 8.99),
    @("Cujo", "Stephen King", "Horror", 8.99),
    @("Blackbirds", "Chuck Wendig", "Horror; Thriller", 8.99),
    @("Hollow Things", "T.S. King", "Horror", 8.99),
    @("Heart Shaped Box", "Joe Hill", "Horror", 8.99),
    @("Buddha Da", "Anne Donovan", "Literature; Scottish", 8.99),
    @("Trainspotting", "Irvine Welsh", "Literature; Scottish", 8.99),
    @("The Crow Road", "Iain Banks", "Literature; Scottish", 8.99),
    @("Klara and the Sun", "Kazuo Ishiguro", "Literature", 8.99),
    @("One: Pot, Pan, Planet", "Anna Jones", "Non-Fiction;Cookery", 8.99),
    @("The Midnight Library", "Matt Haig", "Modern Fiction", 8.99),
    @("Acts of Desperation", "Megan Nolan", "Modern Fiction", 8.99),
    @("Transcendant Kingdom", "Yaa Gyasi", "Modern Fiction", 8.99),
    @("Difficult Women", "Helen Lewis", "Social Sciences", 8.99),
    @("With These Hands", "Pam Ayres", "Biography", 8.99),
    @("Empireland", "Sathnam Sanghera", "History; British History", 8.99),
    @("The Thursday Murder Club", "Richard Osman", "Modern Fiction; Crime; Thriller", 8.99),
    @("Luster", "Raven Leilani", "Modern Fiction", 8.99)
)

$startRow = 7
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $book = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $book[0]
    $ws.Cells.Item($r, 2).Value = $book[1]
    $ws.Cells.Item($r, 3).Value = $book[2]
    $ws.Cells.Item($r, 4).Value = $book[3]
}

# Update selection to reflect the next empty row (matches the diff's activeCell)
$ws.Range("E34").Select() | Out-Null
